$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column headers for K (history) and L (priority)
$ws.Range("K1").Value = "history"
$ws.Range("L1").Value = "priority"

# Row 2 - new columns
$ws.Range("K2").Value = "nw req"
$ws.Range("L2").Value = 1

# Row 4 - Sales Incharge updated first so "Sandeep P" is interned before "Closed"
$ws.Range("I4").Value = "Sandeep P"

# Row 3 - Position Status changed to Closed, new columns
$ws.Range("H3").Value = "Closed"
$ws.Range("K3").Value = "nw req"
$ws.Range("L3").Value = 2

# Row 4 (cont.) - Position Status changed to Closed, new columns
$ws.Range("H4").Value = "Closed"
$ws.Range("K4").Value = "nw req"
$ws.Range("L4").Value = 3

# Row 5 - Position Status changed to Closed, new columns
$ws.Range("H5").Value = "Closed"
$ws.Range("K5").Value = "nw req"
$ws.Range("L5").Value = 1

# Update the active selection to match the saved view state
$ws.Range("C9").Select()
